$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.494.14'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '3.180.42'
$ws.Range("E3").Value = '  -4.21%  '
$ws.Range("D5").Value = "'571.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.67%  '
$ws.Range("D6").Value = "'169.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.02%  '
$ws.Range("D7").Value = "'0.606"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.58%  '
$ws.Range("D9").Value = '3.190.28'
$ws.Range("E9").Value = '  -3.82%  '
$ws.Range("E10").Value = '  -4.01%  '
$ws.Range("D11").Value = "'6.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("E12").Value = '  -3.74%  '
$ws.Range("D13").Value = '3.740.16'
$ws.Range("E13").Value = '  -4.00%  '
$ws.Range("E14").Value = '  -2.27%  '
$ws.Range("D15").Value = '64.532.16'
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("D16").Value = "'25.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.01%  '
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '3.192.71'
$ws.Range("E18").Value = '  -3.84%  '
$ws.Range("D19").Value = "'417.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("E20").Value = '  -2.04%  '
$ws.Range("D21").Value = "'5.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("D22").Value = "'7.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.89%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = "'70.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.97%  '
$ws.Range("D25").Value = "'5.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("E27").Value = '  -4.95%  '
$ws.Range("E28").Value = '  -6.90%  '
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -4.14%  '
$ws.Range("D32").Value = "'21.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("E35").Value = '  -3.70%  '
$ws.Range("E36").Value = '  -3.69%  '
$ws.Range("D37").Value = "'157.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("E38").Value = '  -5.25%  '
$ws.Range("D39").Value = '2.737.86'
$ws.Range("E39").Value = '  -5.06%  '
$ws.Range("E40").Value = '  -5.36%  '
$ws.Range("D41").Value = "'24.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.48%  '
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("D44").Value = "'0.717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.11%  '
$ws.Range("D45").Value = "'5.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("D47").Value = "'0.0264"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.78%  '
$ws.Range("E48").Value = '  -6.46%  '
$ws.Range("D49").Value = "'294.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.99%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").Value = "'2.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -12.78%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.20%  '
